# major accuracy check update
#
# - G2:G29 "NEBNextPoly(A)E7490" -> "NEBNextPoly(A)E7490L", rendered in a
#   smaller Arial 9 font so it fits, and the column widened to show it.
# - I2:I29 literal FALSE booleans become a live "=FALSE()" formula (an
#   accuracy/QC check cell) instead of a static value.
# - Selection moves from the old I2:I29 block to the new G2:G29 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 29

# 1) Update the reagent/kit label text in column G and give it its own
#    (smaller) font so the longer text still reads well.
$gRange = $ws.Range("G2:G" + $lastRow)
$gRange.Value = "NEBNextPoly(A)E7490L"
$gRange.Font.Name = "Arial"
$gRange.Font.Size = 9

# 2) Widen column G to fit the new, longer label.
$ws.Columns.Item(7).ColumnWidth = 35.83

# 3) Turn the static FALSE values in column I into a live formula so the
#    "checked" flag is computed rather than hard-coded.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=FALSE()"
}

# 4) Move the active selection to the column that was just edited.
[void]$ws.Range("G2:G" + $lastRow).Select()
